$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# List of cell updates (Row, Column, NewValue) derived from the diff.
# Columns: B=2 (number_of_inclusions), E=5 (Total_Inclusion_Size), F=6 (Overlap_Area)
$changes = @(
    @{Row=2; Col=2; Value=5}
    @{Row=2; Col=5; Value=1145}
    @{Row=2; Col=6; Value=326}
    @{Row=3; Col=2; Value=163}
    @{Row=3; Col=5; Value=8507}
    @{Row=3; Col=6; Value=468}
    @{Row=4; Col=2; Value=45}
    @{Row=4; Col=5; Value=3031}
    @{Row=4; Col=6; Value=192}
    @{Row=5; Col=2; Value=1}
    @{Row=5; Col=5; Value=21}
    @{Row=6; Col=5; Value=180}
    @{Row=7; Col=2; Value=1}
    @{Row=7; Col=5; Value=16}
    @{Row=8; Col=2; Value=1}
    @{Row=8; Col=5; Value=36}
    @{Row=8; Col=6; Value=0}
    @{Row=9; Col=5; Value=16}
    @{Row=10; Col=2; Value=14}
    @{Row=10; Col=5; Value=881}
    @{Row=11; Col=2; Value=2}
    @{Row=11; Col=5; Value=57}
    @{Row=11; Col=6; Value=0}
    @{Row=12; Col=2; Value=3}
    @{Row=12; Col=5; Value=887}
    @{Row=13; Col=5; Value=275}
    @{Row=14; Col=2; Value=4}
    @{Row=14; Col=5; Value=334}
    @{Row=15; Col=5; Value=240}
    @{Row=16; Col=5; Value=328}
    @{Row=17; Col=2; Value=2}
    @{Row=17; Col=5; Value=100}
    @{Row=19; Col=2; Value=65}
    @{Row=19; Col=5; Value=12040}
    @{Row=20; Col=5; Value=53}
    @{Row=21; Col=2; Value=1}
    @{Row=21; Col=5; Value=12}
    @{Row=22; Col=2; Value=15}
    @{Row=22; Col=5; Value=1921}
    @{Row=22; Col=6; Value=626}
    @{Row=23; Col=2; Value=5}
    @{Row=23; Col=5; Value=1034}
    @{Row=23; Col=6; Value=146}
    @{Row=24; Col=2; Value=5}
    @{Row=24; Col=5; Value=315}
    @{Row=24; Col=6; Value=253}
    @{Row=25; Col=2; Value=30}
    @{Row=25; Col=5; Value=29853}
    @{Row=25; Col=6; Value=379}
    @{Row=26; Col=2; Value=28}
    @{Row=26; Col=5; Value=2652}
    @{Row=26; Col=6; Value=514}
    @{Row=27; Col=2; Value=3}
    @{Row=27; Col=5; Value=67}
    @{Row=28; Col=2; Value=9}
    @{Row=28; Col=5; Value=414}
    @{Row=28; Col=6; Value=26}
    @{Row=29; Col=2; Value=2}
    @{Row=29; Col=5; Value=89}
    @{Row=30; Col=5; Value=521}
    @{Row=31; Col=2; Value=9}
    @{Row=31; Col=5; Value=3324}
    @{Row=31; Col=6; Value=100}
    @{Row=32; Col=5; Value=428}
    @{Row=33; Col=2; Value=3}
    @{Row=33; Col=5; Value=336}
    @{Row=34; Col=2; Value=6}
    @{Row=34; Col=5; Value=1029}
    @{Row=35; Col=5; Value=817}
    @{Row=36; Col=2; Value=14}
    @{Row=36; Col=5; Value=25688}
    @{Row=37; Col=5; Value=1776}
    @{Row=38; Col=2; Value=6}
    @{Row=38; Col=5; Value=2349}
    @{Row=39; Col=2; Value=6}
    @{Row=39; Col=5; Value=1109}
    @{Row=40; Col=2; Value=50}
    @{Row=40; Col=5; Value=9867}
    @{Row=41; Col=5; Value=759}
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, $change.Col).Value = $change.Value
}
